$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 223340.7
$ws.Range("E2").Value = 243581.3
$ws.Range("F2").Value = 250474.2
$ws.Range("G2").Value = 228521.3
$ws.Range("H2").Value = 233350.1
$ws.Range("I2").Value = 264274.2
$ws.Range("J2").Value = 291837.5
$ws.Range("K2").Value = 309894.09999999998
$ws.Range("L2").Value = 303810.8
$ws.Range("M2").Value = 363315.4
$ws.Range("N2").Value = 395723.2
$ws.Range("O2").Value = 396392.1
$ws.Range("P2").Value = 388386.6
$ws.Range("D3").Value = 214181.9
$ws.Range("E3").Value = 233854.3
$ws.Range("F3").Value = 240951.6
$ws.Range("G3").Value = 222052.3
$ws.Range("H3").Value = 226178.7
$ws.Range("I3").Value = 254814.2
$ws.Range("J3").Value = 280621.7
$ws.Range("K3").Value = 293990.7
$ws.Range("L3").Value = 288538.90000000002
$ws.Range("M3").Value = 341445.6
$ws.Range("N3").Value = 371800.8
$ws.Range("O3").Value = 373679.6
$ws.Range("P3").Value = 366237.4
$ws.Range("D4").Value = 11273.8
$ws.Range("E4").Value = 11653.5
$ws.Range("F4").Value = 11783.7
$ws.Range("G4").Value = 10005.700000000001
$ws.Range("H4").Value = 10149.799999999999
$ws.Range("I4").Value = 14389.5
$ws.Range("J4").Value = 17559.599999999999
$ws.Range("K4").Value = 18937.400000000001
$ws.Range("L4").Value = 18855.900000000001
$ws.Range("M4").Value = 24737
$ws.Range("N4").Value = 25995.200000000001
$ws.Range("O4").Value = 25832.7
$ws.Range("P4").Value = 26104.3
$ws.Range("D5").Value = 12667.1
$ws.Range("E5").Value = 13127.9
$ws.Range("F5").Value = 13475.2
$ws.Range("G5").Value = 12266.4
$ws.Range("H5").Value = 12356.8
$ws.Range("I5").Value = 14410.8
$ws.Range("J5").Value = 16984.8
$ws.Range("K5").Value = 18608.7
$ws.Range("L5").Value = 18708.599999999999
$ws.Range("M5").Value = 24185.599999999999
$ws.Range("N5").Value = 25734
$ws.Range("O5").Value = 25733
$ws.Range("P5").Value = 25772.1
$ws.Range("D6").Value = 9203.4
$ws.Range("E6").Value = 9382.7999999999993
$ws.Range("F6").Value = 9704.4
$ws.Range("G6").Value = 8924.6
$ws.Range("H6").Value = 8773.4
$ws.Range("I6").Value = 10331.5
$ws.Range("J6").Value = 12130.9
$ws.Range("K6").Value = 13670.5
$ws.Range("L6").Value = 13999.5
$ws.Range("M6").Value = 18114.8
$ws.Range("N6").Value = 19327.2
$ws.Range("O6").Value = 19526
$ws.Range("P6").Value = 19375
$ws.Range("D7").Value = 1268.5
$ws.Range("E7").Value = 1287.5
$ws.Range("F7").Value = 1455.9
$ws.Range("G7").Value = 1418.6
$ws.Range("H7").Value = 1358.6
$ws.Range("I7").Value = 1342.5
$ws.Range("J7").Value = 1316.3
$ws.Range("K7").Value = 1314
$ws.Range("L7").Value = 1320.4
$ws.Range("M7").Value = 1308.5999999999999
$ws.Range("N7").Value = 1292.8
$ws.Range("O7").Value = 1319.9
$ws.Range("P7").Value = 1296.5
$ws.Range("J8").Value = 129.19999999999999
$ws.Range("K8").Value = 143
$ws.Range("L8").Value = 157
$ws.Range("M8").Value = 172.1
$ws.Range("N8").Value = 189.6
$ws.Range("O8").Value = 200.6
$ws.Range("P8").Value = 209.3
$ws.Range("I9").Value = 394.1
$ws.Range("J9").Value = 422.4
$ws.Range("K9").Value = 436.8
$ws.Range("L9").Value = 424.1
$ws.Range("M9").Value = 438.7
$ws.Range("N9").Value = 448.8
$ws.Range("O9").Value = 454.1
$ws.Range("P9").Value = 455.4
$ws.Range("D10").Value = 9670
$ws.Range("E10").Value = 10732
$ws.Range("F10").Value = 18657.7
$ws.Range("G10").Value = 19555.099999999999
$ws.Range("H10").Value = 20425.8
$ws.Range("I10").Value = 22087
$ws.Range("J10").Value = 24469.4
$ws.Range("K10").Value = 26681.3
$ws.Range("L10").Value = 27368.7
$ws.Range("M10").Value = 30964.5
$ws.Range("N10").Value = 33178.9
$ws.Range("O10").Value = 36823.9
$ws.Range("P10").Value = 37992.1
$ws.Range("D11").Value = 48788.3
$ws.Range("G11").Value = 58270.3
$ws.Range("H11").Value = 61364.2
$ws.Range("I11").Value = 68791.399999999994
$ws.Range("J11").Value = 71527.399999999994
$ws.Range("K11").Value = 74209.399999999994
$ws.Range("L11").Value = 78842.7
$ws.Range("M11").Value = 92875.4
$ws.Range("N11").Value = 100493.1
$ws.Range("O11").Value = 105022.39999999999
$ws.Range("P11").Value = 107015.4
$ws.Range("D12").Value = 3835.3
$ws.Range("E12").Value = 4282.8999999999996
$ws.Range("F12").Value = 4589
$ws.Range("G12").Value = 4935.7
$ws.Range("H12").Value = 5259.2
$ws.Range("I12").Value = 6080
$ws.Range("J12").Value = 6946.1
$ws.Range("K12").Value = 8190.4
$ws.Range("L12").Value = 8462
$ws.Range("M12").Value = 9415.6
$ws.Range("N12").Value = 10281
$ws.Range("O12").Value = 10811.7
$ws.Range("P12").Value = 10872.6
$ws.Range("D13").Value = 18.6267
$ws.Range("E13").Value = 25.511299999999999
$ws.Range("F13").Value = 33.994999999999997
$ws.Range("G13").Value = 41
$ws.Range("H13").Value = 48.8
$ws.Range("I13").Value = 57.5
$ws.Range("J13").Value = 67.5
$ws.Range("K13").Value = 77.2
$ws.Range("L13").Value = 88.5
$ws.Range("M13").Value = 104.4
$ws.Range("N13").Value = 122.9
$ws.Range("O13").Value = 139.80000000000001
$ws.Range("P13").Value = 155.1
$ws.Range("G14").Value = 19745.3
$ws.Range("H14").Value = 19368.5
$ws.Range("I14").Value = 19728.5
$ws.Range("J14").Value = 20692.2
$ws.Range("K14").Value = 20651
$ws.Range("L14").Value = 19191.3
$ws.Range("M14").Value = 22570.400000000001
$ws.Range("N14").Value = 27012
$ws.Range("O14").Value = 25344.2
$ws.Range("P14").Value = 25352.9
$ws.Range("G15").Value = 33.1
$ws.Range("H15").Value = 58
$ws.Range("I15").Value = 68.599999999999994
$ws.Range("J15").Value = 47.2
$ws.Range("L15").Value = 157.5
$ws.Range("M15").Value = 125.9
$ws.Range("N15").Value = 156.4
$ws.Range("O15").Value = 172.3
$ws.Range("P15").Value = 326.5

$ws.Range("Q8").Select() | Out-Null
